# meeting-auto-attender: support for JSON meetings list.
#
# Adds a "Comment (Optional)" column (E) with an example row, and
# corrects the sample Zoom link text shown in B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Comment (Optional)" column header and example comment.
$ws.Range("E2").Value = "Example meeting from Excel"
$ws.Range("E1").Value = "Comment (Optional)"

# Update the example meeting link text (hyperlink target is left as-is).
$ws.Range("B2").Value = "https://us05web.zoom.us/j/87177504375?pwd=jhvL2kxa2ZoQWdicWd1BiS0JLZzadblahblah"

# Leave the active selection on the example link cell.
[void]$ws.Range("B2").Select()
